$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 (latest poll figures)
$ws.Range("A2").Value = 37
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 3
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 5

# Update row 9 (second poll block)
$ws.Range("A9").Value = 44
$ws.Range("B9").Value = 48

# Move selection to D7 to match saved cursor position
$ws.Range("D7").Select()
